$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.886.16"
$ws.Range("E2").Value = "  +0.61%  "
$ws.Range("D3").Value = "2.533.60"
$ws.Range("E3").Value = "  +0.10%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.11"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.53%  "
$ws.Range("E7").Value = "  -0.96%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.535"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.99%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.16"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.99%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0821"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.79%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.60"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.00%  "
$ws.Range("E13").Value = "  -3.67%  "
$ws.Range("D14").Value = "2.924.97"
$ws.Range("E14").Value = "  +0.01%  "
$ws.Range("D15").Value = "2.553.38"
$ws.Range("E15").Value = "  -0.28%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.17"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.21%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.853"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.55%  "
$ws.Range("D18").Value = "42.952.92"
$ws.Range("E18").Value = "  +0.59%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.86"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.93%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.73"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.64%  "
$ws.Range("D21").Value = "0.0₃0969"
$ws.Range("E21").Value = "  -0.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.75"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "253.62"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.38%  "
$ws.Range("E24").Value = "  +0.48%  "
$ws.Range("E25").Value = "  +0.96%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.41"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.12%  "
$ws.Range("E27").Value = "  +0.59%  "
$ws.Range("E28").Value = "  +2.17%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "40.89"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.51%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.51"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.10%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.90"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.23%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "157.67"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.25%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.16"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.95%  "
$ws.Range("B34").Value = "WEMIXToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.72"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.59%  "
$ws.Range("B35").Value = "Celestia"
$ws.Range("C35").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "19.23"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.49%  "
$ws.Range("B36").Value = "LidoDAOToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.33"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.95%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0793"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.38%  "
$ws.Range("E38").Value = "  +0.13%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.47"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +9.45%  "
$ws.Range("E40").Value = "  -0.99%  "
$ws.Range("E41").Value = "  -12.27%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.82"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.62%  "
$ws.Range("E43").Value = "  +0.64%  "
$ws.Range("E44").Value = "  +0.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.28"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.47%  "
$ws.Range("D46").Value = "2.011.63"
$ws.Range("E46").Value = "  -1.63%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.19"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.28%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "84.25"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.69%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "106.52"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.65%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "75.70"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.08%  "
$ws.Range("D51").Value = "2.781.94"
$ws.Range("E51").Value = "  +0.15%  "
